# Login Tests and PIM Tests added.

$wb = $excel.ActiveWorkbook

# --- Rename existing sheets ---
$wsValid = $wb.Worksheets.Item("testValidAdminLogin")
$wsInvalid = $wb.Worksheets.Item("Sheet2")
$wsInvalid.Name = "testInvalidAdminLogin"
$wsPim = $wb.Worksheets.Item("Sheet3")
$wsPim.Name = "testPIMPageLinks"

# --- testInvalidAdminLogin (Sheet2) data ---
$wsInvalid.Range("A1").Value = "Username"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "invalidUserName"
$wsInvalid.Range("B2").Value = "invalidPassword"

# --- testPIMPageLinks (Sheet3) data ---
$wsPim.Range("A1").Value = "Username"
$wsPim.Range("B1").Value = "Password"
$wsPim.Range("C1").Value = "TabNames"
$wsPim.Range("A2").Value = "Admin"
$wsPim.Range("B2").Value = "admin123"
$wsPim.Range("C2").Value = "Employee List;Add Employee;Reports"

# --- New sheet testCreateUser (added after the last existing sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCreate = $wb.Worksheets.Add($null, $lastSheet)
$wsCreate.Name = "testCreateUser"
$wsCreate.Range("A1").Value = "Username"
$wsCreate.Range("B1").Value = "Password"
$wsCreate.Range("C1").Value = "FirstName"
$wsCreate.Range("D1").Value = "MiddleName"
$wsCreate.Range("E1").Value = "LastName"
$wsCreate.Range("A2").Value = "Admin"
$wsCreate.Range("B2").Value = "admin123"
$wsCreate.Range("C2").Value = "Harry"
$wsCreate.Range("D2").Value = "James"
$wsCreate.Range("E2").Value = "Potter"

# --- Selections on each sheet (order matters: last Select/Activate wins the active tab) ---
$wsValid.Range("A1:B2").Select()
$wsInvalid.Range("E5").Select()
$wsPim.Range("A1:B2").Select()

# testCreateUser ends up active, matching activeTab="3" / tabSelected="1"
$wsCreate.Range("M13").Select()
$wsCreate.Activate()
